$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename header row (row 1) cells: "_old" -> "_FV2410", "_new" -> "_FV2504"
for ($c = 1; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Value2
    if ($val -ne $null) {
        if ($val -like "*_old") {
            $cell.Value = $val -replace "_old$", "_FV2410"
        } elseif ($val -like "*_new") {
            $cell.Value = $val -replace "_new$", "_FV2504"
        }
    }
}

# Create a table (ListObject) over A1:U64 with headers
$range = $ws.Range("A1:U64")
$table = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $range, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$table.Name = "Table1"
$table.TableStyle = ""

# Freeze the first row (pane split) and set selection in the frozen view
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

$wb.Save()
